$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so the holdings cells can be
# refreshed, then restore protection with the same password afterwards.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclaimer banner.
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-28 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns with the latest
# model holdings figures.
$ws.Range("D2").Value = 0.03502798557723331
$ws.Range("E2").Value = 0.0003999999999999559
$ws.Range("D3").Value = 0.02037115551618042
$ws.Range("E3").Value = 0.0003901677721420693
$ws.Range("D4").Value = 0.01932220109135994
$ws.Range("E4").Value = 0.006012024048096309
$ws.Range("D5").Value = 0.03803198562033684
$ws.Range("E5").Value = 0.006269592476489061
$ws.Range("D6").Value = 0.03425003312373691
$ws.Range("E6").Value = 0.001600640256102626
$ws.Range("D7").Value = 0.01979590683085715
$ws.Range("E7").Value = -0.00009651578033009756
$ws.Range("D8").Value = 0.0369037023612709
$ws.Range("E8").Value = 0.005349882302589259
$ws.Range("D9").Value = 0.02047509310247505
$ws.Range("E9").Value = 0.002060378034578481
$ws.Range("D10").Value = 0.02547744609148439
$ws.Range("E10").Value = 0.002099790020997805
$ws.Range("D11").Value = 0.02411836025108774
$ws.Range("E11").Value = -0.0002640612622127891
$ws.Range("D12").Value = 0.05761122843858717
$ws.Range("E12").Value = -0.00118063754427411
$ws.Range("D13").Value = 0.02489916589322752
$ws.Range("E13").Value = 0.003683241252302016
$ws.Range("D14").Value = 0.02659640553766213
$ws.Range("E14").Value = 0.004084197298146464
$ws.Range("D15").Value = 0.03212119774719358
$ws.Range("E15").Value = 0.003553028957186122
$ws.Range("D16").Value = 0.01910133372048386
$ws.Range("E16").Value = 0.002069857697283295
$ws.Range("D17").Value = 0.03173662867790346
$ws.Range("E17").Value = 0.006421576497030124
$ws.Range("D18").Value = 0.04210032692701623
$ws.Range("E18").Value = 0.002529317084387195
$ws.Range("D19").Value = 0.1257403292126241
$ws.Range("E19").Value = 0.003992015968063978
$ws.Range("D20").Value = 0.008896649788442432
$ws.Range("E20").Value = -0.0003951527924130893
$ws.Range("D21").Value = 0.01515685137080557
$ws.Range("E21").Value = 0.00529436679373152
$ws.Range("D22").Value = 0.01767545269595375
$ws.Range("E22").Value = 0.01406380163669319
$ws.Range("D23").Value = 0.01517223820956095
$ws.Range("E23").Value = 0.003294289897511016
$ws.Range("D24").Value = 0.02179887073340985
$ws.Range("E24").Value = 0.006934673366834287
$ws.Range("D25").Value = 0.0127089683641755
$ws.Range("E25").Value = -0.005372011818425837
$ws.Range("D26").Value = 0.04227814171925556
$ws.Range("E26").Value = 0.002331887201735494
$ws.Range("D27").Value = 0.02391033222910687
$ws.Range("E27").Value = -0.00009801999607916834
$ws.Range("D28").Value = 0.04579673471411191
$ws.Range("E28").Value = 0.001421800947867258
$ws.Range("D29").Value = 0.05643179158312916
$ws.Range("E29").Value = -0.002477437621659972
$ws.Range("D30").Value = 0.01354072380351853
$ws.Range("E30").Value = -0.008744534665833914
$ws.Range("D31").Value = 0.02058015158430814
$ws.Range("E31").Value = 0.001918649270913342
$ws.Range("D32").Value = 0.01341253411375515
$ws.Range("E32").Value = 0.005223171889838563
$ws.Range("D33").Value = 0.04178887081671277
$ws.Range("E33").Value = 0
$ws.Range("D34").Value = 0.0171712025530332
$ws.Range("E34").Value = 0.002471648735097309
$ws.Range("E35").Value = 0.002327768859723411

# Restore sheet protection to match the original workbook state.
$ws.Protect("D382")
